# "añadido el diagrama de bd"
# Insert a new blank row above the existing row 5 (the user-story table),
# pushing the rest of the sheet down by one row, then move the active
# selection to the new C5 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Rows("5:5").Insert()

$ws.Range("C5").Select() | Out-Null
